$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21 (ALC)
$ws.Range("H21").Value = 13949.5
$ws.Range("I21").Value = 8831.666999999999
$ws.Range("J21").Value = 16142.857
$ws.Range("K21").Value = 8831.666999999999
$ws.Range("L21").Value = 16142.857
$ws.Range("M21").Value = -8363.666999999999
$ws.Range("N21").Value = -17078.857

# Row 23 (ALC)
$ws.Range("H23").Value = 13949.5
$ws.Range("I23").Value = 8831.666999999999
$ws.Range("J23").Value = 16142.857
$ws.Range("K23").Value = 8831.666999999999
$ws.Range("L23").Value = 16142.857
$ws.Range("M23").Value = -8597.666999999999
$ws.Range("N23").Value = -16610.857

# Row 40 (ALC)
$ws.Range("H40").Value = 3500
$ws.Range("J40").Value = 1000
$ws.Range("L40").Value = 1000
$ws.Range("N40").Value = -1350

# Row 62 (ALC)
$ws.Range("H62").Value = 2158
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 2816
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 2816
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -4064

# Row 65 (ALC)
$ws.Range("H65").Value = 2158
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 2816
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 14080
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -20320

# Row 134 (ALC)
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0

# Row 136 (ALC)
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("N136").Value = 0

# Row 137 (ALC)
$ws.Range("H137").Value = 999
$ws.Range("I137").Value = 999
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 2997
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -447

# Row 141 (ALC)
$ws.Range("H141").Value = 2166.6667
$ws.Range("J141").Value = 5000
$ws.Range("L141").Value = 15000
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
# Row 74 (ARM)
$ws.Range("H74").Value = 1484.1666
$ws.Range("I74").Value = 1484.1666
$ws.Range("K74").Value = 1484.1666
$ws.Range("M74").Value = -610.1666

# Row 77 (ARM)
$ws.Range("H77").Value = 1484.1666
$ws.Range("I77").Value = 1484.1666
$ws.Range("K77").Value = 7420.833000000001
$ws.Range("M77").Value = -3052.833000000001

# Row 88 (ARM)
$ws.Range("H88").Value = 2471.1428
$ws.Range("I88").Value = 2449.5
$ws.Range("K88").Value = 2449.5
$ws.Range("M88").Value = -2043.5

# Row 91 (ARM)
$ws.Range("H91").Value = 2471.1428
$ws.Range("I91").Value = 2449.5
$ws.Range("K91").Value = 2449.5
$ws.Range("M91").Value = -1045.5

# Row 139 (ARM)
$ws.Range("H139").Value = 79749.25
$ws.Range("J139").Value = 79749.25
$ws.Range("L139").Value = 79749.25
$ws.Range("N139").Value = -90029.25

$ws = $wb.Worksheets.Item("BSM")
# Row 38 (BSM)
$ws.Range("H38").Value = 31000
$ws.Range("J38").Value = 31000
$ws.Range("L38").Value = 31000
$ws.Range("N38").Value = -31832

# Row 86 (BSM)
$ws.Range("H86").Value = 2176.6667
$ws.Range("I86").Value = 1716.8334
$ws.Range("J86").Value = 4935.6665
$ws.Range("K86").Value = 1716.8334
$ws.Range("L86").Value = 4935.6665
$ws.Range("M86").Value = -593.8334
$ws.Range("N86").Value = -7181.6665

# Row 88 (BSM)
$ws.Range("H88").Value = 27999
$ws.Range("J88").Value = 27999
$ws.Range("L88").Value = 27999
$ws.Range("N88").Value = -28811

# Row 89 (BSM)
$ws.Range("H89").Value = 2176.6667
$ws.Range("I89").Value = 1716.8334
$ws.Range("J89").Value = 4935.6665
$ws.Range("K89").Value = 8584.166999999999
$ws.Range("L89").Value = 24678.3325
$ws.Range("M89").Value = -2968.166999999999
$ws.Range("N89").Value = -35910.3325

# Row 91 (BSM)
$ws.Range("H91").Value = 27999
$ws.Range("J91").Value = 27999
$ws.Range("L91").Value = 27999
$ws.Range("N91").Value = -30807

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 7689.769
$ws.Range("I31").Value = 6332.8335
$ws.Range("J31").Value = 8852.857
$ws.Range("K31").Value = 6332.8335
$ws.Range("L31").Value = 8852.857
$ws.Range("M31").Value = -6037.8335
$ws.Range("N31").Value = -9442.857

# Row 34 (CRP)
$ws.Range("H34").Value = 7689.769
$ws.Range("I34").Value = 6332.8335
$ws.Range("J34").Value = 8852.857
$ws.Range("K34").Value = 6332.8335
$ws.Range("L34").Value = 8852.857
$ws.Range("M34").Value = -6130.8335
$ws.Range("N34").Value = -9256.857

# Row 51 (CRP)
$ws.Range("H51").Value = 23363.334
$ws.Range("I51").Value = 23363.334
$ws.Range("K51").Value = 23363.334
$ws.Range("M51").Value = -22627.334

# Row 61 (CRP)
$ws.Range("H61").Value = 23363.334
$ws.Range("I61").Value = 23363.334
$ws.Range("K61").Value = 23363.334
$ws.Range("M61").Value = -23015.334

# Row 62 (CRP)
$ws.Range("H62").Value = 3375
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -2748

# Row 65 (CRP)
$ws.Range("H65").Value = 3375
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -13740

$ws = $wb.Worksheets.Item("CUL")
# Row 21 (CUL)
$ws.Range("H21").Value = 298.5
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 298.5
$ws.Range("K21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("M21").Value = 895.5
$ws.Range("N21").Value = -1241.5

# Row 86 (CUL)
$ws.Range("H86").Value = 774.25
$ws.Range("I86").Value = 198
$ws.Range("K86").Value = 594
$ws.Range("M86").Value = 592

# Row 89 (CUL)
$ws.Range("H89").Value = 774.25
$ws.Range("I89").Value = 198
$ws.Range("K89").Value = 1782
$ws.Range("M89").Value = 4146

$ws = $wb.Worksheets.Item("GSM")
# Row 4 (GSM)
$ws.Range("H4").Value = 250
$ws.Range("J4").Value = 250
$ws.Range("L4").Value = 250
$ws.Range("N4").Value = -474

# Row 70 (GSM)
$ws.Range("H70").Value = 166669660
$ws.Range("I70").Value = 4497.5
$ws.Range("K70").Value = 4497.5
$ws.Range("M70").Value = -4227.5

# Row 73 (GSM)
$ws.Range("H73").Value = 166669660
$ws.Range("I73").Value = 4497.5
$ws.Range("K73").Value = 4497.5
$ws.Range("M73").Value = -3561.5

# Row 80 (GSM)
$ws.Range("H80").Value = 28250
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 83 (GSM)
$ws.Range("H83").Value = 28250
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Row 132 (GSM)
$ws.Range("H132").Value = 112
$ws.Range("I132").Value = 112
$ws.Range("K132").Value = 336
$ws.Range("M132").Value = 2194

$ws = $wb.Worksheets.Item("LTW")
# Row 122 (LTW)
$ws.Range("H122").Value = 5071.2856
$ws.Range("I122").Value = 4899.8
$ws.Range("K122").Value = 14699.4
$ws.Range("M122").Value = -12249.4

$ws = $wb.Worksheets.Item("WVR")
# Row 11 (WVR)
$ws.Range("H11").Value = 487.5
$ws.Range("I11").Value = 450
$ws.Range("K11").Value = 450
$ws.Range("M11").Value = -308

# Row 43 (WVR)
$ws.Range("H43").Value = 9318.75
$ws.Range("I43").Value = 758.3333
$ws.Range("K43").Value = 758.3333
$ws.Range("M43").Value = -609.3333

# Row 103 (WVR)
$ws.Range("H103").Value = 17473
$ws.Range("J103").Value = 17473
$ws.Range("L103").Value = 17473
$ws.Range("N103").Value = -19817

# Row 122 (WVR)
$ws.Range("H122").Value = 1602.75
$ws.Range("I122").Value = 1514.0714
$ws.Range("J122").Value = 2223.5
$ws.Range("K122").Value = 4542.2142
$ws.Range("L122").Value = 6670.5
$ws.Range("M122").Value = -2092.2142
$ws.Range("N122").Value = -11570.5
